$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-05-14 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-15 Wednesday", 2)

# Update the division-problem table. Each cell is addressed by its
# explicit (row, column) position so that values which coincidentally
# match other cells' old/new text are never confused with one another.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="63÷5="},
    @{Row=1;  Col=2; New="81÷6="},
    @{Row=1;  Col=3; New="42÷5="},
    @{Row=1;  Col=4; New="58÷5="},
    @{Row=1;  Col=5; New="31÷9="},

    @{Row=5;  Col=1; New="66÷7="},
    @{Row=5;  Col=2; New="71÷5="},
    @{Row=5;  Col=3; New="16÷7="},
    @{Row=5;  Col=4; New="88÷5="},
    @{Row=5;  Col=5; New="44÷7="},

    @{Row=9;  Col=1; New="26÷8="},
    @{Row=9;  Col=2; New="68÷5="},
    @{Row=9;  Col=3; New="42÷4="},
    @{Row=9;  Col=4; New="17÷5="},
    @{Row=9;  Col=5; New="79÷4="},

    @{Row=13; Col=1; New="46÷3="},
    @{Row=13; Col=2; New="87÷2="},
    @{Row=13; Col=3; New="87÷4="},
    @{Row=13; Col=4; New="46÷7="},
    @{Row=13; Col=5; New="78÷2="},

    @{Row=17; Col=1; New="50÷8="},
    @{Row=17; Col=2; New="89÷6="},
    @{Row=17; Col=3; New="46÷9="},
    @{Row=17; Col=4; New="36÷7="},
    @{Row=17; Col=5; New="57÷2="}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text is replaced, leaving the cell's formatting untouched.
    $rng.End = $rng.End - 1
    $rng.Text = $u.New
}
